# "Code optimization on DNC algos"
#
# Fills in the empirical trial measurements for the "Naive DNC" (rows 8-12)
# and refreshes the "Enhanced DNC" (rows 14-18) sections, plus appends a new
# n=1,000,000 data point (row 19) to the "Enhanced DNC" table. The M (Average)
# and N ("(n, avg)" label) columns already carry AVERAGE/concat formulas for
# rows 8-12 (they currently show #DIV/0! because C:L are empty) and for rows
# 14-18 (already populated) - so once the trial columns are filled those
# formulas recompute on their own. Row 19 is brand new, so its M/N formulas
# are copied down (format + formula) from row 18 and then re-pointed at row 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Naive DNC (rows 8-12): trial columns C:L were empty -> AVERAGE()/"(...)"
# formulas in M/N errored with #DIV/0!. Fill in the 10 trial values per row.
# ---------------------------------------------------------------------

# n = 10
$ws.Cells.Item(8, 3).Value = 0.0000758171
$ws.Cells.Item(8, 4).Value = 0.0000846386
$ws.Cells.Item(8, 5).Value = 0.0000629425
$ws.Cells.Item(8, 6).Value = 0.0000658035
$ws.Cells.Item(8, 7).Value = 0.0000479221
$ws.Cells.Item(8, 8).Value = 0.0000674725
$ws.Cells.Item(8, 9).Value = 0.0000588894
$ws.Cells.Item(8, 10).Value = 0.0000479221
$ws.Cells.Item(8, 11).Value = 0.0000481606
$ws.Cells.Item(8, 12).Value = 0.0000441074

# n = 100
$ws.Cells.Item(9, 3).Value = 0.0003962517
$ws.Cells.Item(9, 4).Value = 0.0003905296
$ws.Cells.Item(9, 5).Value = 0.0003988743
$ws.Cells.Item(9, 6).Value = 0.0003976822
$ws.Cells.Item(9, 7).Value = 0.0004045963
$ws.Cells.Item(9, 8).Value = 0.0004787445
$ws.Cells.Item(9, 9).Value = 0.0003905296
$ws.Cells.Item(9, 10).Value = 0.0003900528
$ws.Cells.Item(9, 11).Value = 0.0003917217
$ws.Cells.Item(9, 12).Value = 0.0003886223

# n = 1000
$ws.Cells.Item(10, 3).Value = 0.0039718151
$ws.Cells.Item(10, 4).Value = 0.0040888786
$ws.Cells.Item(10, 5).Value = 0.0041782856
$ws.Cells.Item(10, 6).Value = 0.0038986206
$ws.Cells.Item(10, 7).Value = 0.0038070679
$ws.Cells.Item(10, 8).Value = 0.0047163963
$ws.Cells.Item(10, 9).Value = 0.003872633
$ws.Cells.Item(10, 10).Value = 0.0037899017
$ws.Cells.Item(10, 11).Value = 0.0037956238
$ws.Cells.Item(10, 12).Value = 0.004267931

# n = 10000
$ws.Cells.Item(11, 3).Value = 0.0472183228
$ws.Cells.Item(11, 4).Value = 0.0450963974
$ws.Cells.Item(11, 5).Value = 0.0452141762
$ws.Cells.Item(11, 6).Value = 0.044169426
$ws.Cells.Item(11, 7).Value = 0.0440707207
$ws.Cells.Item(11, 8).Value = 0.0443511009
$ws.Cells.Item(11, 9).Value = 0.0443711281
$ws.Cells.Item(11, 10).Value = 0.0441226959
$ws.Cells.Item(11, 11).Value = 0.044062376
$ws.Cells.Item(11, 12).Value = 0.0492372513

# n = 100000
$ws.Cells.Item(12, 3).Value = 0.5273954868
$ws.Cells.Item(12, 4).Value = 0.5369346142
$ws.Cells.Item(12, 5).Value = 0.5256202221
$ws.Cells.Item(12, 6).Value = 0.5384941101
$ws.Cells.Item(12, 7).Value = 0.5620825291
$ws.Cells.Item(12, 8).Value = 0.5287666321
$ws.Cells.Item(12, 9).Value = 0.5310969353
$ws.Cells.Item(12, 10).Value = 0.5618369579
$ws.Cells.Item(12, 11).Value = 0.5255272388
$ws.Cells.Item(12, 12).Value = 0.533427

# ---------------------------------------------------------------------
# Enhanced DNC (rows 14-18): re-measured trial values (replace old ones).
# ---------------------------------------------------------------------

# n = 10
$ws.Cells.Item(14, 3).Value = 0.0000693798
$ws.Cells.Item(14, 4).Value = 0.0000550747
$ws.Cells.Item(14, 5).Value = 0.0000436306
$ws.Cells.Item(14, 6).Value = 0.0000450611
$ws.Cells.Item(14, 7).Value = 0.0000367165
$ws.Cells.Item(14, 8).Value = 0.0000443459
$ws.Cells.Item(14, 9).Value = 0.0000350475
$ws.Cells.Item(14, 10).Value = 0.0000388622
$ws.Cells.Item(14, 11).Value = 0.0000340939
$ws.Cells.Item(14, 12).Value = 0.0000333786

# n = 100
$ws.Cells.Item(15, 3).Value = 0.0003170967
$ws.Cells.Item(15, 4).Value = 0.0003116131
$ws.Cells.Item(15, 5).Value = 0.0003144741
$ws.Cells.Item(15, 6).Value = 0.0003049374
$ws.Cells.Item(15, 7).Value = 0.0003147125
$ws.Cells.Item(15, 8).Value = 0.000305891
$ws.Cells.Item(15, 9).Value = 0.0003020763
$ws.Cells.Item(15, 10).Value = 0.0003089905
$ws.Cells.Item(15, 11).Value = 0.0002975464
$ws.Cells.Item(15, 12).Value = 0.0003061295

# n = 1000
$ws.Cells.Item(16, 3).Value = 0.0032274723
$ws.Cells.Item(16, 4).Value = 0.0031468868
$ws.Cells.Item(16, 5).Value = 0.0032169819
$ws.Cells.Item(16, 6).Value = 0.0031776428
$ws.Cells.Item(16, 7).Value = 0.0031192303
$ws.Cells.Item(16, 8).Value = 0.0031518936
$ws.Cells.Item(16, 9).Value = 0.0031163692
$ws.Cells.Item(16, 10).Value = 0.0031385422
$ws.Cells.Item(16, 11).Value = 0.0031421185
$ws.Cells.Item(16, 12).Value = 0.0031671524

# n = 10000
$ws.Cells.Item(17, 3).Value = 0.0406410694
$ws.Cells.Item(17, 4).Value = 0.0405409336
$ws.Cells.Item(17, 5).Value = 0.0405299664
$ws.Cells.Item(17, 6).Value = 0.0406548977
$ws.Cells.Item(17, 7).Value = 0.0405220985
$ws.Cells.Item(17, 8).Value = 0.0406138897
$ws.Cells.Item(17, 9).Value = 0.0405135155
$ws.Cells.Item(17, 10).Value = 0.0411868095
$ws.Cells.Item(17, 11).Value = 0.0406839848
$ws.Cells.Item(17, 12).Value = 0.0405249596

# n = 100000
$ws.Cells.Item(18, 3).Value = 0.5544199944
$ws.Cells.Item(18, 4).Value = 0.542617321
$ws.Cells.Item(18, 5).Value = 0.5580203533
$ws.Cells.Item(18, 6).Value = 0.5660443306
$ws.Cells.Item(18, 7).Value = 0.5544099808
$ws.Cells.Item(18, 8).Value = 0.5558817387
$ws.Cells.Item(18, 9).Value = 0.5537266731
$ws.Cells.Item(18, 10).Value = 0.5552794933
$ws.Cells.Item(18, 11).Value = 0.5548655987
$ws.Cells.Item(18, 12).Value = 0.5699870586

# ---------------------------------------------------------------------
# Enhanced DNC: new row 19 (n = 1,000,000). Copy row 18's M/N cell formats
# (number format + font + borders) down into row 19 first, same as a
# fill-down/autofill would, then write B19 and the trial data C19:L19,
# and finally point the copied M19/N19 formulas at row 19.
# ---------------------------------------------------------------------

$ws.Range("M18").Copy()
$ws.Range("M19").PasteSpecial(-4122)
$ws.Range("N18").Copy()
$ws.Range("N19").PasteSpecial(-4122)

$ws.Cells.Item(19, 2).Value = 1000000

$ws.Cells.Item(19, 3).Value = 9.766307354
$ws.Cells.Item(19, 4).Value = 9.6560795307
$ws.Cells.Item(19, 5).Value = 10.7907721996
$ws.Cells.Item(19, 6).Value = 10.0300133228
$ws.Cells.Item(19, 7).Value = 9.7327427864
$ws.Cells.Item(19, 8).Value = 10.7889122963
$ws.Cells.Item(19, 9).Value = 11.4387202263
$ws.Cells.Item(19, 10).Value = 11.5064268112
$ws.Cells.Item(19, 11).Value = 10.70480299
$ws.Cells.Item(19, 12).Value = 10.8377747536

$ws.Range("M19").Formula = "=AVERAGE(C19:L19)"
$ws.Range("N19").Formula = "=""("" & B19 & "", "" & M19 & "")"""

# ---------------------------------------------------------------------
# Cursor/selection moved from L13 to J14 as part of this edit.
# ---------------------------------------------------------------------
$ws.Range("J14").Select() | Out-Null
